# Add a "Precio" (Price) column H to the Components sheet and fill with
# sequential values 1..113 for rows 2..114. Also update the selection /
# scroll position to match the recorded view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column H
$ws.Range("H1").Value = "Precio"

# Fill H2:H114 with sequential numbers 1..113
for ($row = 2; $row -le 114; $row++) {
    $ws.Cells.Item($row, 8).Value = $row - 1
}

# Update view state: scroll + selection
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("H2:H114").Select()
